# Daily attendance processing - 2026-01-01 08:40:10
# Applies the daily refresh of the "Session Analysis Results" sheet:
#   1. Swap the "Recorded By" name order from "System, <email>" to
#      "<email>, System" for every session row that was (re)processed.
#   2. Refresh the workbook-level "Missing Sessions" / "Pending Sessions"
#      counters (K7:L8 metric block).
#   3. Refresh the per-group "Missing" / "Pending" counters (P21:Q26) in
#      the Group Statistics table.
#   4. Six sessions dated 01/01/2026 have now fallen past due without
#      attendance being captured, so they flip from the "Pending" look
#      (yellow row style) to the "Not Recorded" look (green row style),
#      matching the style/status already used for earlier missed
#      sessions (e.g. row 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. "Recorded By" (column G) name-order swap -------------------------
$recordedByRows = @(8,9,10,12,14,15,17,18,34,35,36,38,40,41,43,44,60,61,62,64,66,67,69,70,86,87,88,90,92,93,95,96,112,113,114,116,118,119,121,122,138,139,140,142,144,145,147,148,164,167,170,174,191,194,197,201,218,221,224,228,245,248,251,255,272,275,278,282,299,302,305,309)

foreach ($r in $recordedByRows) {
    $ws.Range("G$r").Value = "dnasr281@gmail.com, System"
}

# --- 2. Workbook-level metric block (K7:L8) -------------------------------
$ws.Range("L7").Value = 9    # Missing Sessions
$ws.Range("L8").Value = 108  # Pending Sessions

# --- 3. Group Statistics table: Missing (P) / Pending (Q) columns --------
$ws.Range("P21").Value = 1
$ws.Range("Q21").Value = 9

$ws.Range("P22").Value = 1
$ws.Range("Q22").Value = 9

$ws.Range("P23").Value = 1
$ws.Range("Q23").Value = 9

$ws.Range("P24").Value = 2
$ws.Range("Q24").Value = 9

$ws.Range("P25").Value = 1
$ws.Range("Q25").Value = 9

$ws.Range("P26").Value = 1
$ws.Range("Q26").Value = 9

# --- 4. Rows that flipped from "Pending" to "Not Recorded" ---------------
# Copy the exact row style already used for "Not Recorded" sessions (row 3)
# onto each newly-overdue row, then update the Status text in column I.
$notRecordedRows = @(175,202,229,256,283,310)

$ws.Range("A3:I3").Copy()
foreach ($r in $notRecordedRows) {
    $ws.Range("A$r`:I$r").PasteSpecial(-4122)
    $ws.Range("I$r").Value = "Not Recorded"
}
